$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every touched cell in the sheet is plain text (inline strings in the
# original file - prices such as "0.05280" or "1.000" are NOT numbers).
# Force text format first, cell by cell, so Excel does not "helpfully"
# reinterpret a numeric-looking string as a number (which would eat
# significant trailing/leading zeros, e.g. "0.05280" -> 0.0528) when we
# write the new values below.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.713.50"
$ws.Range("E2").Value = "  -0.50%  "

# Row 3
$ws.Range("D3").Value = "1.868.31"
$ws.Range("E3").Value = "  -0.13%  "

# Row 4
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "0.7314"
$ws.Range("E5").Value = "  -0.44%  "

# Row 6
$ws.Range("D6").Value = "240.93"
$ws.Range("E6").Value = "  -0.42%  "

# Row 7
$ws.Range("D7").Value = "0.9988"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "0.3132"
$ws.Range("E8").Value = "  -0.64%  "

# Row 9
$ws.Range("D9").Value = "0.07101"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").Value = "24.43"

# Row 11
$ws.Range("D11").Value = "0.08235"
$ws.Range("E11").Value = "  -1.91%  "

# Row 12
$ws.Range("D12").Value = "0.7464"
$ws.Range("E12").Value = "  -0.73%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.314"
$ws.Range("E13").Value = "  -1.74%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.828.24"
$ws.Range("E14").Value = "  -1.51%  "

# Row 15
$ws.Range("D15").Value = "92.47"
$ws.Range("E15").Value = "  -0.06%  "

# Row 16
$ws.Range("D16").Value = "29.669.43"
$ws.Range("E16").Value = "  -0.65%  "

# Row 17
$ws.Range("D17").Value = "6.022"
$ws.Range("E17").Value = "  -0.36%  "

# Row 18
$ws.Range("D18").Value = "248.63"
$ws.Range("E18").Value = "  +2.34%  "

# Row 19
$ws.Range("D19").Value = "13.37"
$ws.Range("E19").Value = "  -1.44%  "

# Row 20
$ws.Range("D20").Value = "0.000007806"
$ws.Range("E20").Value = "  -0.17%  "

# Row 21
$ws.Range("D21").Value = "0.9959"
$ws.Range("E21").Value = "  -0.28%  "

# Row 22
$ws.Range("D22").Value = "2.129.17"
$ws.Range("E22").Value = "  +1.07%  "

# Row 23
$ws.Range("D23").Value = "0.9997"
$ws.Range("E23").Value = "  -0.13%  "

# Row 24
$ws.Range("D24").Value = "7.714"
$ws.Range("E24").Value = "  -2.53%  "

# Row 25
$ws.Range("E25").Value = "  -1.98%  "

# Row 26
$ws.Range("D26").Value = "9.168"
$ws.Range("E26").Value = "  -1.64%  "

# Row 27
$ws.Range("D27").Value = "162.81"
$ws.Range("E27").Value = "  -0.76%  "

# Row 28
$ws.Range("D28").Value = "18.54"
$ws.Range("E28").Value = "  -0.16%  "

# Row 29
$ws.Range("D29").Value = "2.015"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("D30").Value = "1.434"
$ws.Range("E30").Value = "  -2.56%  "

# Row 31
$ws.Range("D31").Value = "4.539"
$ws.Range("E31").Value = "  -1.86%  "

# Row 32
$ws.Range("D32").Value = "1.522"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("D33").Value = "4.192"
$ws.Range("E33").Value = "  -2.38%  "

# Row 34
$ws.Range("D34").Value = "0.05280"
$ws.Range("E34").Value = "  -0.97%  "

# Row 35
$ws.Range("E35").Value = "  -0.11%  "

# Row 36
$ws.Range("D36").Value = "0.7540"
$ws.Range("E36").Value = "  +0.27%  "

# Row 37
$ws.Range("D37").Value = "0.9976"
$ws.Range("E37").Value = "  -0.22%  "

# Row 38
$ws.Range("D38").Value = "2.694"
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("D39").Value = "0.01930"
$ws.Range("E39").Value = "  -1.04%  "

# Row 40
$ws.Range("E40").Value = "  -0.45%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("D42").Value = "5.993"
$ws.Range("E42").Value = "  -1.34%  "

# Row 43
$ws.Range("D43").Value = "0.8647"
$ws.Range("E43").Value = "  +0.37%  "

# Row 44
$ws.Range("D44").Value = "71.13"

# Row 45
$ws.Range("D45").Value = "1.048.47"
$ws.Range("E45").Value = "  -5.08%  "

# Row 46
$ws.Range("D46").Value = "103.91"
$ws.Range("E46").Value = "  +1.07%  "

# Row 47
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  -0.10%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.481"
$ws.Range("E48").Value = "  -2.97%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.820"
$ws.Range("E49").Value = "  -1.12%  "

# Row 50
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "2.899"
$ws.Range("E50").Value = "  -5.15%  "

# Row 51
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "1.989.58"
$ws.Range("E51").Value = "  -1.29%  "
